# 20250106 calculo de fuste mixto
#
# Rename the "Tipo Cálculo" options from full words to their abbreviations
# ("Drenado" -> "d", "No drenado" -> "nd") and flip row 4 (H4) from the
# "drenado" case to the "no drenado" case.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# H2 keeps meaning "Drenado", just abbreviated.
$ws.Range("H2").Value = "d"

# H3 keeps meaning "No drenado", just abbreviated.
$ws.Range("H3").Value = "nd"

# H4 switches from "Drenado" to "No drenado" (mixed shaft calculation).
$ws.Range("H4").Value = "nd"
